$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Modifica Antigua BD": the old "Periodo Mora" value 2508 becomes 2509 for
# every worker row that carried it (E19:E23 all shared that value).
$ws.Range("E19:E23").Value = "2509"

# Center the "Periodo Mora" column values (E16:E23) horizontally, matching
# the alignment already used elsewhere in the table.
$ws.Range("E16:E23").HorizontalAlignment = -4108
